$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10337
$ws1.Range("F4").Value = 2526
$ws1.Range("F9").Value = 774
$ws1.Range("F12").Value = 1086
$ws1.Range("F13").Value = 3223
$ws1.Range("F14").Value = 2400
$ws1.Range("F16").Value = 2152
$ws1.Range("F21").Value = 575
$ws1.Range("F22").Value = 63
$ws1.Range("F23").Value = 248
$ws1.Range("F28").Value = 377
$ws1.Range("F29").Value = 6
$ws1.Range("F31").Value = 389
$ws1.Range("F32").Value = 599
$ws1.Range("F35").Value = 258
$ws1.Range("F38").Value = 464
$ws1.Range("F39").Value = 447
$ws1.Range("F40").Value = 1712
$ws1.Range("F42").Value = 443
$ws1.Range("F43").Value = 51
$ws1.Range("F44").Value = 455
$ws1.Range("F45").Value = 1021
$ws1.Range("F47").Value = 364

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10337
$ws4.Range("F9").Value = 774
$ws4.Range("F10").Value = 1086
$ws4.Range("F11").Value = 3223
$ws4.Range("F12").Value = 2400
$ws4.Range("F13").Value = 2152
$ws4.Range("F16").Value = 575
$ws4.Range("F17").Value = 63
$ws4.Range("F18").Value = 248
$ws4.Range("F23").Value = 377
$ws4.Range("F24").Value = 6
$ws4.Range("F26").Value = 389
$ws4.Range("F27").Value = 599
$ws4.Range("F33").Value = 258
$ws4.Range("F36").Value = 465
$ws4.Range("F38").Value = 447
$ws4.Range("F39").Value = 1712
$ws4.Range("F44").Value = 443
$ws4.Range("F45").Value = 51
$ws4.Range("F46").Value = 455
$ws4.Range("F47").Value = 1021
$ws4.Range("F48").Value = 364
